$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D price cells keep their original text formatting
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '29.135.22'
$ws.Range('E2').Value = '  +8.64%  '
$ws.Range('D3').Value = '1.822.41'
$ws.Range('E3').Value = '  +5.72%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '246.18'
$ws.Range('E5').Value = '  +2.54%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').Value = '0.4927'
$ws.Range('E7').Value = '  +3.12%  '
$ws.Range('D8').Value = '44.43'
$ws.Range('E8').Value = '  +7.85%  '
$ws.Range('D9').Value = '0.2777'
$ws.Range('E9').Value = '  +7.97%  '
$ws.Range('D10').Value = '0.06382'
$ws.Range('E10').Value = '  +4.02%  '
$ws.Range('D11').Value = '1.822.28'
$ws.Range('E11').Value = '  +5.67%  '
$ws.Range('D12').Value = '16.65'
$ws.Range('E12').Value = '  +4.73%  '
$ws.Range('D13').Value = '0.07112'
$ws.Range('E13').Value = '  +3.31%  '
$ws.Range('D14').Value = '0.6448'
$ws.Range('E14').Value = '  +7.71%  '
$ws.Range('D15').Value = '84.06'
$ws.Range('E15').Value = '  +9.65%  '
$ws.Range('D16').Value = '4.701'
$ws.Range('E16').Value = '  +6.40%  '
$ws.Range('D17').Value = '29.133.77'
$ws.Range('E17').Value = '  +9.40%  '
$ws.Range('D18').Value = '0.9996'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').Value = '0.000007308'
$ws.Range('E19').Value = '  +3.50%  '
$ws.Range('D20').Value = '0.9999'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').Value = '12.22'
$ws.Range('E21').Value = '  +8.54%  '
$ws.Range('D22').Value = '2.056.31'
$ws.Range('E22').Value = '  +5.60%  '
$ws.Range('D23').Value = '4.545'
$ws.Range('E23').Value = '  +3.86%  '
$ws.Range('D24').Value = '8.843'
$ws.Range('E24').Value = '  +6.13%  '
$ws.Range('D25').Value = '5.369'
$ws.Range('E25').Value = '  +6.84%  '
$ws.Range('D26').Value = '143.63'
$ws.Range('E26').Value = '  +2.47%  '
$ws.Range('D27').Value = '130.38'
$ws.Range('E27').Value = '  +22.82%  '
$ws.Range('D28').Value = '16.40'
$ws.Range('E28').Value = '  +8.26%  '
$ws.Range('D29').Value = '1.884'
$ws.Range('E29').Value = '  +5.97%  '
$ws.Range('D30').Value = '1.401'
$ws.Range('E30').Value = '  +1.32%  '
$ws.Range('D31').Value = '4.120'
$ws.Range('E31').Value = '  +5.05%  '
$ws.Range('D32').Value = '0.08328'
$ws.Range('E32').Value = '  +5.78%  '
$ws.Range('D33').Value = '3.775'
$ws.Range('E33').Value = '  +4.19%  '
$ws.Range('D34').Value = '0.04928'
$ws.Range('E34').Value = '  +8.52%  '
$ws.Range('D35').Value = '1.095'
$ws.Range('E35').Value = '  +10.23%  '
$ws.Range('D36').Value = '2.696'
$ws.Range('E36').Value = '  +3.80%  '
$ws.Range('D37').Value = '0.6686'
$ws.Range('E37').Value = '  +9.30%  '
$ws.Range('D38').Value = '2.289'
$ws.Range('E38').Value = '  +16.29%  '
$ws.Range('D39').Value = '2.680'
$ws.Range('E39').Value = '  +7.63%  '
$ws.Range('D40').Value = '0.9501'
$ws.Range('E40').Value = '  +3.56%  '
$ws.Range('D41').Value = '6.153'
$ws.Range('E41').Value = '  +7.44%  '
$ws.Range('D42').Value = '0.01589'
$ws.Range('E42').Value = '  +7.39%  '
$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').Value = '101.91'
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('D45').Value = '0.4059'
$ws.Range('E45').Value = '  +6.92%  '
$ws.Range('D46').Value = '7.153'
$ws.Range('E46').Value = '  +6.58%  '
$ws.Range('D47').Value = '0.1217'
$ws.Range('E47').Value = '  +6.54%  '
$ws.Range('D48').Value = '0.05547'
$ws.Range('E48').Value = '  +3.63%  '
$ws.Range('D49').Value = '31.72'
$ws.Range('E49').Value = '  +6.67%  '
$ws.Range('D50').Value = '8.077'
$ws.Range('E50').Value = '  +2.86%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.301'
$ws.Range('E51').Value = '  +5.93%  '
